# Weekly data refresh ("Fruta / hortaliza, semanal"): insert 3 new price
# records at the top of the data table (rows 49-51), pushing the existing
# records down by three rows (old 49..109 -> new 52..112).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right above the first data row (row 49), shifting the
# remaining 61 records (old rows 49-109) down to rows 52-112.
$ws.Rows("49:51").Insert()

# --- New record 1 (row 49) ---
$ws.Range("A49").Value2 = 10
$ws.Range("B49").Value2 = "Vega Modelo de Temuco"
$ws.Range("C49").Value2 = "La Araucanía"
$ws.Range("D49").Value2 = 44966
$ws.Range("E49").Value2 = 9
$ws.Range("F49").Value2 = 100112030
$ws.Range("G49").Value2 = "Poroto granado"
$ws.Range("H49").Value2 = "Sin especificar"
$ws.Range("I49").Value2 = "Primera"
$ws.Range("J49").Value2 = 80
$ws.Range("K49").Value2 = 45000
$ws.Range("L49").Value2 = 45000
$ws.Range("M49").Value2 = 45000
$ws.Range("N49").Value2 = "$/saco 25 kilos"
$ws.Range("O49").Value2 = "Región de La Araucanía"
$ws.Range("P49").Value2 = 1800
$ws.Range("Q49").Value2 = 25
$ws.Range("R49").Value2 = "Hortaliza"

# --- New record 2 (row 50) ---
$ws.Range("A50").Value2 = 10
$ws.Range("B50").Value2 = "Vega Modelo de Temuco"
$ws.Range("C50").Value2 = "La Araucanía"
$ws.Range("D50").Value2 = 44966
$ws.Range("E50").Value2 = 9
$ws.Range("F50").Value2 = 100112030
$ws.Range("G50").Value2 = "Poroto granado"
$ws.Range("H50").Value2 = "Sin especificar"
$ws.Range("I50").Value2 = "Primera"
$ws.Range("J50").Value2 = 30
$ws.Range("K50").Value2 = 45000
$ws.Range("L50").Value2 = 45000
$ws.Range("M50").Value2 = 45000
$ws.Range("N50").Value2 = "$/saco 25 kilos"
$ws.Range("O50").Value2 = "Región del Maule"
$ws.Range("P50").Value2 = 1800
$ws.Range("Q50").Value2 = 25
$ws.Range("R50").Value2 = "Hortaliza"

# --- New record 3 (row 51) ---
$ws.Range("A51").Value2 = 10
$ws.Range("B51").Value2 = "Vega Modelo de Temuco"
$ws.Range("C51").Value2 = "La Araucanía"
$ws.Range("D51").Value2 = 44966
$ws.Range("E51").Value2 = 9
$ws.Range("F51").Value2 = 100112030
$ws.Range("G51").Value2 = "Poroto granado"
$ws.Range("H51").Value2 = "Sin especificar"
$ws.Range("I51").Value2 = "Segunda"
$ws.Range("J51").Value2 = 50
$ws.Range("K51").Value2 = 40000
$ws.Range("L51").Value2 = 40000
$ws.Range("M51").Value2 = 40000
$ws.Range("N51").Value2 = "$/saco 25 kilos"
$ws.Range("O51").Value2 = "Región del Maule"
$ws.Range("P51").Value2 = 1600
$ws.Range("Q51").Value2 = 25
$ws.Range("R51").Value2 = "Hortaliza"
